$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Cell="D2"; Value="69.373.38"},
    @{Cell="E2"; Value="  +1.72%  "},
    @{Cell="D3"; Value="3.939.05"},
    @{Cell="E3"; Value="  +0.31%  "},
    @{Cell="E4"; Value="  -0.03%  "},
    @{Cell="D5"; Value="491.71"},
    @{Cell="E5"; Value="  +1.07%  "},
    @{Cell="D6"; Value="147.14"},
    @{Cell="E6"; Value="  +0.41%  "},
    @{Cell="E8"; Value="  +0.02%  "},
    @{Cell="D9"; Value="0.736"},
    @{Cell="E9"; Value="  +0.49%  "},
    @{Cell="E10"; Value="  +4.09%  "},
    @{Cell="E11"; Value="  -3.53%  "},
    @{Cell="E12"; Value="  +0.57%  "},
    @{Cell="D13"; Value="10.49"},
    @{Cell="E13"; Value="  -1.78%  "},
    @{Cell="D14"; Value="4.570.37"},
    @{Cell="E14"; Value="  +0.39%  "},
    @{Cell="D15"; Value="3.946.98"},
    @{Cell="E15"; Value="  +0.66%  "},
    @{Cell="E16"; Value="  -4.45%  "},
    @{Cell="E17"; Value="  -0.68%  "},
    @{Cell="E18"; Value="  -1.20%  "},
    @{Cell="E19"; Value="  +2.29%  "},
    @{Cell="D20"; Value="69.383.74"},
    @{Cell="E20"; Value="  +1.53%  "},
    @{Cell="D21"; Value="440.02"},
    @{Cell="E21"; Value="  -1.78%  "},
    @{Cell="E22"; Value="  +1.74%  "},
    @{Cell="D23"; Value="14.57"},
    @{Cell="E23"; Value="  -1.65%  "},
    @{Cell="D24"; Value="89.48"},
    @{Cell="E24"; Value="  +0.90%  "},
    @{Cell="D25"; Value="12.09"},
    @{Cell="E25"; Value="  +10.55%  "},
    @{Cell="E26"; Value="  +3.32%  "},
    @{Cell="D27"; Value="11.16"},
    @{Cell="E27"; Value="  -3.74%  "},
    @{Cell="D28"; Value="37.33"},
    @{Cell="E28"; Value="  -4.65%  "},
    @{Cell="E29"; Value="  -3.68%  "},
    @{Cell="D30"; Value="707.84"},
    @{Cell="E30"; Value="  +2.57%  "},
    @{Cell="E31"; Value="  +0.03%  "},
    @{Cell="E32"; Value="  +0.05%  "},
    @{Cell="E33"; Value="  +0.98%  "},
    @{Cell="D34"; Value="0.473"},
    @{Cell="E34"; Value="  +26.58%  "},
    @{Cell="D35"; Value="0.0₃0905"},
    @{Cell="E35"; Value="  -4.75%  "},
    @{Cell="D36"; Value="61.55"},
    @{Cell="E36"; Value="  +4.61%  "},
    @{Cell="D37"; Value="6.07"},
    @{Cell="E37"; Value="  +4.29%  "},
    @{Cell="D38"; Value="40.82"},
    @{Cell="E38"; Value="  -2.50%  "},
    @{Cell="E39"; Value="  -0.04%  "},
    @{Cell="D40"; Value="0.998"},
    @{Cell="E40"; Value="  -0.16%  "},
    @{Cell="E41"; Value="  +0.13%  "},
    @{Cell="E42"; Value="  +2.12%  "},
    @{Cell="E43"; Value="  +4.84%  "},
    @{Cell="E44"; Value="  -0.84%  "},
    @{Cell="D45"; Value="3.01"},
    @{Cell="E45"; Value="  +2.04%  "},
    @{Cell="B47"; Value="BabyDogeCoin"},
    @{Cell="C47"; Value="https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"},
    @{Cell="D47"; Value="0.0₆0366"},
    @{Cell="E47"; Value="  +10.90%  "},
    @{Cell="B48"; Value="ApeXProtocol"},
    @{Cell="C48"; Value="https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"},
    @{Cell="D48"; Value="3.35"},
    @{Cell="E48"; Value="  +6.78%  "},
    @{Cell="D49"; Value="3.06"},
    @{Cell="E49"; Value="  +7.57%  "},
    @{Cell="D50"; Value="3.38"},
    @{Cell="E50"; Value="  -1.49%  "},
    @{Cell="E51"; Value="  -3.28%  "}
)

foreach ($item in $changes) {
    $cell = $ws.Range($item.Cell)
    # Force text interpretation so numeric-looking strings (e.g. "10.49") are not
    # auto-converted into actual numbers, then restore the default cell style so no
    # extraneous style index is introduced, matching the original formatting.
    $cell.NumberFormat = "@"
    $cell.Value = $item.Value
    $cell.Style = "Normal"
}

Write-Output ("Applied {0} cell changes" -f $changes.Count)
